$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to be stored as text so numeric-looking
# strings (e.g. "577.06", "0.0000237") are not coerced to numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.322.37'
$ws.Range('E2').Value = '  +2.76%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.249.16'
$ws.Range('E3').Value = '  +5.23%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.06'
$ws.Range('E5').Value = '  +2.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.51'
$ws.Range('E6').Value = '  +6.91%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.241.11'
$ws.Range('E8').Value = '  +5.41%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  +3.62%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.10'
$ws.Range('E10').Value = '  +9.44%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.166'
$ws.Range('E11').Value = '  +4.32%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.493'
$ws.Range('E12').Value = '  +4.06%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.96'
$ws.Range('E13').Value = '  +3.87%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000237'
$ws.Range('E14').Value = '  +4.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.765.12'
$ws.Range('E15').Value = '  +5.37%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '561.00'
$ws.Range('E16').Value = '  +12.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.429.01'
$ws.Range('E17').Value = '  +2.93%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.248.68'
$ws.Range('E18').Value = '  +5.33%  '

$ws.Range('E19').Value = '  +2.87%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.15'
$ws.Range('E20').Value = '  +5.36%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.49'
$ws.Range('E21').Value = '  +3.87%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.748'
$ws.Range('E22').Value = '  +6.94%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.89'
$ws.Range('E23').Value = '  +8.36%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.62'
$ws.Range('E24').Value = '  +5.88%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.40'
$ws.Range('E25').Value = '  +3.28%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.39'
$ws.Range('E27').Value = '  +17.35%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.97'
$ws.Range('E28').Value = '  +5.94%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').Value = '  +6.22%  '

$ws.Range('E30').Value = '  +5.52%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.80'
$ws.Range('E31').Value = '  +2.49%  '

$ws.Range('E32').Value = '  -0.14%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  +4.11%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '570.68'
$ws.Range('E34').Value = '  +9.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.78'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.43'
$ws.Range('E36').Value = '  +5.57%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0461'
$ws.Range('E37').Value = '  +12.90%  '

$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.73'
$ws.Range('E38').Value = '  +3.55%  '

$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.11'
$ws.Range('E39').Value = '  +15.33%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0869'
$ws.Range('E40').Value = '  +6.85%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.129'
$ws.Range('E41').Value = '  +5.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.152.60'
$ws.Range('E42').Value = '  +6.41%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.65'
$ws.Range('E43').Value = '  +1.80%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.276'
$ws.Range('E44').Value = '  +10.20%  '

$ws.Range('E45').Value = '  +6.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.76'
$ws.Range('E46').Value = '  +4.11%  '

$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₃0563'
$ws.Range('E47').Value = '  +2.92%  '

$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.09%  '

$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.63'
$ws.Range('E49').Value = '  +3.26%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.114'
$ws.Range('E50').Value = '  +2.99%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.28'
$ws.Range('E51').Value = '  +8.30%  '
